$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 23:58"

# --- Rows whose ranking swapped (country label + stats move together) ---
# Egipto <-> Republica Dominicana (rows 34/35)
$ws.Range("A34").Value = "Republica Dominicana"
$ws.Range("B34").Value = 99898
$ws.Range("C34").Value = 565
$ws.Range("D34").Value = 73228
$ws.Range("E34").Value = 24806
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 19
$ws.Range("H34").Value = 1864

$ws.Range("A35").Value = "Egipto"
$ws.Range("B35").Value = 99863
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 78108
$ws.Range("E35").Value = 16225
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 5530

# Mauritania <-> Zimbabue (rows 105/106)
$ws.Range("A105").Value = "Zimbabue"
$ws.Range("B105").Value = 7298
$ws.Range("C105").Value = 461
$ws.Range("D105").Value = 5455
$ws.Range("E105").Value = 1633
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 4
$ws.Range("H105").Value = 210

$ws.Range("A106").Value = "Mauritania"
$ws.Range("B106").Value = 7142
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 6669
$ws.Range("E106").Value = 313
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 160

# --- Plain stat refresh rows (same country, updated counters) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6480924
$ws.Range("C4").Value = 20674
$ws.Range("D4").Value = 3741341
$ws.Range("E4").Value = 2546107
$ws.Range("G4").Value = 228
$ws.Range("H4").Value = 193476

# Row 5: India
$ws.Range("B5").Value = 4277584
$ws.Range("C5").Value = 75022
$ws.Range("D5").Value = 3321420
$ws.Range("E5").Value = 883348
$ws.Range("G5").Value = 1129
$ws.Range("H5").Value = 72816

# Row 6: Brasil
$ws.Range("B6").Value = 4147794
$ws.Range("C6").Value = 10188
$ws.Range("E6").Value = 703607
$ws.Range("G6").Value = 274
$ws.Range("H6").Value = 126960

# Row 23: Irak
$ws.Range("B23").Value = 264684
$ws.Range("C23").Value = 4314
$ws.Range("D23").Value = 202859
$ws.Range("E23").Value = 54236
$ws.Range("G23").Value = 77
$ws.Range("H23").Value = 7589

# Row 82: Costa de Marfil
$ws.Range("B82").Value = 18701
$ws.Range("C82").Value = 113
$ws.Range("D82").Value = 17599
$ws.Range("E82").Value = 983

# Row 89: Zambia
$ws.Range("B89").Value = 12836
$ws.Range("C89").Value = 60
$ws.Range("D89").Value = 11748
$ws.Range("E89").Value = 793

# Row 95: Guinea
$ws.Range("B95").Value = 9816
$ws.Range("C95").Value = 18
$ws.Range("D95").Value = 8956
$ws.Range("E95").Value = 798
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 62

# Row 103: Haiti
$ws.Range("B103").Value = 8362
$ws.Range("C103").Value = 2
$ws.Range("E103").Value = 2215

# Row 108: Malaui
$ws.Range("B108").Value = 5621
$ws.Range("C108").Value = 7
$ws.Range("D108").Value = 3590
$ws.Range("E108").Value = 1855
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 176

# Row 114: Suazilandia
$ws.Range("B114").Value = 4884
$ws.Range("C114").Value = 31
$ws.Range("D114").Value = 4029
$ws.Range("E114").Value = 761

# Row 120: Ruanda
$ws.Range("B120").Value = 4409
$ws.Range("C120").Value = 35
$ws.Range("D120").Value = 2271
$ws.Range("E120").Value = 2119

# Row 142: Trinidad yTobago
$ws.Range("B142").Value = 2254
$ws.Range("C142").Value = 4
$ws.Range("D142").Value = 734
$ws.Range("E142").Value = 1485
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 35

# Row 155: Togo
$ws.Range("B155").Value = 1493
$ws.Range("C155").Value = 5
$ws.Range("D155").Value = 1114
$ws.Range("E155").Value = 346
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 33

# Row 166: Republica del Chad
$ws.Range("B166").Value = 1040
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 919

Write-Output "edits applied"
